$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F ("Packs Unboxed") before the existing "Expiration Date" column.
$ws.Range("F1").EntireColumn.Insert()

# Header for the new column (reuse the bold header formatting already used by the row)
$ws.Cells.Item(1, 6).Value = "Packs Unboxed"
$ws.Cells.Item(1, 6).Font.Bold = $true

# Fill "Packs Unboxed" = 0 for existing data rows (2-6), matching the wrapped-text
# formatting already used by the other data cells in those rows.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 6).WrapText = $true
}

# Append new batch rows (7-11)
$newRows = @(
    @("(Aciclovir) 200 mg Tablet 5x: ₱5.0", "1", 5, 0, 5, 0, 44239),
    @("Mucosolve (Ambroxol) 15 mg/60 ml Syrup 1x: ₱9.0", "45", 7, 0, 7, 0, 44246),
    @("(Aciclovir) 200 mg Tablet 1x: ₱2.0", "4", 5, 0, 5, 0, 44242),
    @("(Allopurinol) 100 mg Tablet 10x: ₱2.0", "456", 9, 0, 9, 0, 44244),
    @("Cisflem (Carbo) 125 mg/60 ml Syrup 50x: ₱9.0", "5678", 6, 0, 6, 0, 44251)
)

$row = 7
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 1).WrapText = $true

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 2).WrapText = $true

    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 3).WrapText = $true

    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 4).WrapText = $true

    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 5).WrapText = $true

    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 6).WrapText = $true

    $ws.Cells.Item($row, 7).Value = $data[6]
    $ws.Cells.Item($row, 7).WrapText = $true
    $ws.Cells.Item($row, 7).NumberFormat = "yyyy-mm-dd"

    $row++
}
